# "Fruta / hortaliza, semanal" weekly refresh: the daily Haba price-series
# rows (Fecha / Volumen / Precio minimo / Precio maximo / Precio promedio
# ponderado / Precio $/Kg -> columns D, J, K, L, M, P) get re-dated and
# re-shuffled across the existing rows 2-23 (row 13 keeps its own data).
# All other columns (Mercado, Region, Categoria, Unidad, Origen, etc.)
# stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot original D, J, K, L, M, P values for rows 2-23 before any writes,
# since the edit permutes these values among rows.
$origD = @{}
$origJ = @{}
$origK = @{}
$origL = @{}
$origM = @{}
$origP = @{}
for ($r = 2; $r -le 23; $r++) {
    $origD[$r] = $ws.Cells.Item($r, 4).Value()
    $origJ[$r] = $ws.Cells.Item($r, 10).Value()
    $origK[$r] = $ws.Cells.Item($r, 11).Value()
    $origL[$r] = $ws.Cells.Item($r, 12).Value()
    $origM[$r] = $ws.Cells.Item($r, 13).Value()
    $origP[$r] = $ws.Cells.Item($r, 16).Value()
}

# Row -> source-row map describing where the new Fecha/Volumen/Precio block comes from.
$sourceRow = @{}
$sourceRow[2] = 14
$sourceRow[3] = 6
$sourceRow[4] = 8
$sourceRow[5] = 7
$sourceRow[6] = 11
$sourceRow[7] = 22
$sourceRow[8] = 5
$sourceRow[9] = 17
$sourceRow[10] = 21
$sourceRow[11] = 2
$sourceRow[12] = 23
$sourceRow[13] = 13
$sourceRow[14] = 10
$sourceRow[15] = 16
$sourceRow[16] = 12
$sourceRow[17] = 15
$sourceRow[18] = 3
$sourceRow[19] = 20
$sourceRow[20] = 19
$sourceRow[21] = 18
$sourceRow[22] = 4
$sourceRow[23] = 9

foreach ($r in $sourceRow.Keys) {
    $src = $sourceRow[$r]
    if ($src -eq $r) { continue }
    $ws.Cells.Item($r, 4).Value = $origD[$src]
    $ws.Cells.Item($r, 10).Value = $origJ[$src]
    $ws.Cells.Item($r, 11).Value = $origK[$src]
    $ws.Cells.Item($r, 12).Value = $origL[$src]
    $ws.Cells.Item($r, 13).Value = $origM[$src]
    $ws.Cells.Item($r, 16).Value = $origP[$src]
}

Write-Host "Done"
